$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '30.463.29'
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 5).Value = '  -0.58%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.896.74'
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).Value = '  +0.44%  '
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.9993'
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(4, 5).Value = '  -0.33%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '238.35'
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = '  +0.72%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.9988'
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = '  -0.26%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.4902'
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).Value = '  +0.59%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.2937'
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value = '  +0.70%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.06712'
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = '  +0.51%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '1.900.90'
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = '  +0.69%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '17.01'
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = '  +2.02%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.07319'
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value = '  +0.94%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '5.144'
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = '  +2.75%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '87.73'
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = '  -1.59%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.6664'
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).Value = '  +0.70%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '30.432.45'
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).Value = '  -0.49%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '13.47'
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(17, 5).Value = '  +3.58%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.000007849'
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).Value = '  -0.46%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '0.9993'
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = '  -0.24%  '
$ws.Cells.Item(20, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '2.125.80'
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = '  -0.39%  '
$ws.Cells.Item(21, 2).Value = 'Uniswap'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '5.352'
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = '  +12.89%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '0.9990'
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = '  -0.34%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '191.03'
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value = '  -0.04%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '6.137'
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = '  +0.89%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '9.486'
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = '  +1.95%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '163.07'
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value = '  +2.38%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '18.26'
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = '  -0.32%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '1.936'
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value = '  +5.93%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '1.456'
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = '  +3.63%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '4.369'
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = '  +2.82%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.09198'
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = '  +2.27%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '4.044'
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).Value = '  +2.87%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.05189'
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = '  +0.27%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.7415'
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = '  +1.57%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '1.103'
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value = '  +2.04%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.711'
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value = '  +0.59%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.01808'
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = '  -0.43%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '2.672'
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).Value = '  +0.52%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.9222'
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = '  +0.18%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '2.048'
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value = '  +0.57%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.4394'
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value = '  +0.03%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '5.936'
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = '  +3.85%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '106.38'
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = '  +1.96%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.9937'
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = '  -0.55%  '
$ws.Cells.Item(45, 2).Value = 'Aave'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '68.85'
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = '  +20.35%  '
$ws.Cells.Item(46, 2).Value = 'Algorand'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.1370'
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = '  +2.80%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '7.598'
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = '  +3.62%  '
$ws.Cells.Item(48, 2).Value = 'Elrond'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '35.04'
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value = '  +5.62%  '
$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '8.958'
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = '  +3.82%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.05825'
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = '  -0.10%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.3939'
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = '  -3.55%  '
